$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = 6
$ws.Range("F10").Value = -2
